# #3473 swapped out two properties
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BPS Data")

# Two record IDs were swapped out for refreshed values, and the L3 score was
# recalculated.
$ws.Range("L3").Value = 227440.2
$ws.Range("B6").Value = 22482006
$ws.Range("B10").Value = 22482007

# The active selection moved from the whole data block to column L (rows 2-10),
# anchored at L2.
$ws.Range("L2:L10").Select()
$ws.Application.ActiveCell = $ws.Range("L2")

# Touch the sheet's page setup (portrait) so a pageSetup element is persisted.
$ws.PageSetup.Orientation = 1
